$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data refresh reorders the existing records (rows 2-13) into a new
# chronological arrangement. The underlying set of records is unchanged; only
# the row each record occupies has shifted. Columns D (Fecha), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado),
# O (Origen) and P (Precio $/Kg) are rewritten per row to reflect the new order.

$data = @(
    @{ Row = 2;  D = 44243; J = 80; K = 10000; L = 11000; M = 10375; O = "Provincia de Quillota"; P = 173 },
    @{ Row = 3;  D = 44179; J = 15; K = 7000;  L = 7000;  M = 7000;  O = "Provincia de Limarí";   P = 117 },
    @{ Row = 4;  D = 44333; J = 25; K = 10000; L = 11000; M = 10400; O = "Provincia de Limarí";   P = 173 },
    @{ Row = 5;  D = 44277; J = 25; K = 10000; L = 10000; M = 10000; O = "Provincia de Limarí";   P = 167 },
    @{ Row = 6;  D = 44186; J = 15; K = 7000;  L = 7000;  M = 7000;  O = "Provincia de Limarí";   P = 117 },
    @{ Row = 7;  D = 44585; J = 30; K = 11000; L = 11000; M = 11000; O = "Provincia de Limarí";   P = 183 },
    @{ Row = 8;  D = 44312; J = 30; K = 10000; L = 10000; M = 10000; O = "Provincia de Limarí";   P = 167 },
    @{ Row = 9;  D = 44291; J = 20; K = 9000;  L = 9000;  M = 9000;  O = "Provincia de Limarí";   P = 150 },
    @{ Row = 10; D = 44405; J = 45; K = 9000;  L = 9000;  M = 9000;  O = "Provincia de Quillota"; P = 180 },
    @{ Row = 11; D = 44315; J = 25; K = 10000; L = 10000; M = 10000; O = "Provincia de Limarí";   P = 167 },
    @{ Row = 12; D = 44284; J = 35; K = 10000; L = 10000; M = 10000; O = "Provincia de Limarí";   P = 167 },
    @{ Row = 13; D = 44200; J = 10; K = 9000;  L = 9000;  M = 9000;  O = "Provincia de Limarí";   P = 150 }
)

foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Cells.Item($r, 4).Value = $rec.D    # D: Fecha
    $ws.Cells.Item($r, 10).Value = $rec.J   # J: Volumen
    $ws.Cells.Item($r, 11).Value = $rec.K   # K: Precio mínimo
    $ws.Cells.Item($r, 12).Value = $rec.L   # L: Precio máximo
    $ws.Cells.Item($r, 13).Value = $rec.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 15).Value = $rec.O   # O: Origen
    $ws.Cells.Item($r, 16).Value = $rec.P   # P: Precio $/Kg
}
